# rerouted to free up all adc inputs for use.
#
# Moves the ADC-labelled pin assignments (L column, "A5-PWM"/"A5-LED"/
# "A7-PWM"/"A8-PWM"/"A8-LED") off the dedicated-ADC GPIOs (rows 29-32,
# which previously held "free adc"/"A7-PWM"/"A8-PWM"/"A8-LED") and onto
# earlier, non-ADC GPIOs (rows 14-18), clearing the L column on rows
# 29-32 entirely (and the leftover "free" note in M23) so every ADC
# input is free for use.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reroute analog labels onto rows 14-18 (previously unused GPIOs, "-").
$ws.Range("L14").Value = "A5-PWM"
$ws.Range("L15").Value = "A5-LED"
$ws.Range("L16").Value = "A8-PWM"
$ws.Range("L17").Value = "A7-PWM"
$ws.Range("L18").Value = "A8-LED"

# Clear the stray "free" annotation in M23.
$ws.Range("M23").Clear()

# Free up the dedicated ADC GPIOs (rows 29-32) entirely.
$ws.Range("L29").Clear()
$ws.Range("L30").Clear()
$ws.Range("L31").Clear()
$ws.Range("L32").Clear()

# Update the remembered selection to match the author's final cursor spot.
$ws.Range("L31").Select()
